# Update loading_percent values for Case_3_86 (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.04903736828209
$ws.Range("C2").Value = 9.635483778350785
$ws.Range("D2").Value = 5.283040473236317
$ws.Range("F2").Value = 32.58119619088582
$ws.Range("G2").Value = 3.630382603296448
$ws.Range("I2").Value = 24.67712829652954
$ws.Range("M2").Value = 19.97068109752786
$ws.Range("B3").Value = 11.61348796866855
$ws.Range("C3").Value = 9.088424819342666
$ws.Range("D3").Value = 5.282301887135214
$ws.Range("F3").Value = 31.96636925061955
$ws.Range("G3").Value = 3.634338372838222
$ws.Range("I3").Value = 24.43862110990388
$ws.Range("M3").Value = 19.37051244932731
$ws.Range("B4").Value = 11.34203739984712
$ws.Range("C4").Value = 8.737879024933424
$ws.Range("D4").Value = 5.282339168516243
$ws.Range("F4").Value = 31.59232899849037
$ws.Range("G4").Value = 3.636889321778836
$ws.Range("I4").Value = 24.29670738514301
$ws.Range("M4").Value = 19.0009097405942
$ws.Range("B5").Value = 11.23062014375998
$ws.Range("C5").Value = 8.591507609217587
$ws.Range("D5").Value = 5.282476263061575
$ws.Range("F5").Value = 31.44097248317535
$ws.Range("G5").Value = 3.637959686365243
$ws.Range("I5").Value = 24.24005681321824
$ws.Range("M5").Value = 18.85026831684099
$ws.Range("B6").Value = 11.21207769328428
$ws.Range("C6").Value = 8.56699546292819
$ws.Range("D6").Value = 5.282506345485499
$ws.Range("F6").Value = 31.41591007732244
$ws.Range("G6").Value = 3.638139285631195
$ws.Range("I6").Value = 24.23072247026235
$ws.Range("M6").Value = 18.82526044323104
$ws.Range("B7").Value = 11.3405377289375
$ws.Range("C7").Value = 8.735919027567389
$ws.Range("D7").Value = 5.28234052585245
$ws.Range("F7").Value = 31.59028316843727
$ws.Range("G7").Value = 3.636903632075962
$ws.Range("I7").Value = 24.29593854392802
$ws.Range("M7").Value = 18.99887790047422
$ws.Range("B8").Value = 11.8998163627169
$ws.Range("C8").Value = 9.449986357910285
$ws.Range("D8").Value = 5.28268304769168
$ws.Range("F8").Value = 32.36860229989759
$ws.Range("G8").Value = 3.631721292199265
$ws.Range("I8").Value = 24.59397620335673
$ws.Range("M8").Value = 19.76411697961651
$ws.Range("B9").Value = 12.95625085061589
$ws.Range("C9").Value = 10.72835657531639
$ws.Range("D9").Value = 5.287315702895419
$ws.Range("F9").Value = 33.91379922458731
$ws.Range("G9").Value = 3.622521456838308
$ws.Range("I9").Value = 25.21244472375709
$ws.Range("M9").Value = 21.24572212415582
$ws.Range("B10").Value = 13.69788813672482
$ws.Range("C10").Value = 11.59452545477743
$ws.Range("D10").Value = 5.293224876311932
$ws.Range("F10").Value = 35.04929556525845
$ws.Range("G10").Value = 3.61634079351384
$ws.Range("I10").Value = 25.68483642728819
$ws.Range("M10").Value = 22.30992548471249
$ws.Range("B11").Value = 14.02603709867504
$ws.Range("C11").Value = 12.00285603937058
$ws.Range("D11").Value = 5.29647615096927
$ws.Range("F11").Value = 35.56362182883305
$ws.Range("G11").Value = 3.613652841519106
$ws.Range("I11").Value = 25.90303090203815
$ws.Range("M11").Value = 22.78641877336087
$ws.Range("B12").Value = 14.14884926060178
$ws.Range("C12").Value = 12.1539622573804
$ws.Range("D12").Value = 5.297789797138251
$ws.Range("F12").Value = 35.75788326216966
$ws.Range("G12").Value = 3.612652625751135
$ws.Range("I12").Value = 25.98607636737838
$ws.Range("M12").Value = 22.96558600019056
$ws.Range("B13").Value = 14.12246579855879
$ws.Range("C13").Value = 12.12157510568244
$ws.Range("D13").Value = 5.297503189319714
$ws.Range("F13").Value = 35.71607078490672
$ws.Range("G13").Value = 3.612867257035157
$ws.Range("I13").Value = 25.96817324121529
$ws.Range("M13").Value = 22.92705823658411
$ws.Range("B14").Value = 14.03617066006642
$ws.Range("C14").Value = 12.01535808424609
$ws.Range("D14").Value = 5.296582564471177
$ws.Range("F14").Value = 35.57961493325161
$ws.Range("G14").Value = 3.613570200108333
$ws.Range("I14").Value = 25.90985500179023
$ws.Range("M14").Value = 22.8011853708182
$ws.Range("B15").Value = 13.9831200485906
$ws.Range("C15").Value = 11.94983921234527
$ws.Range("D15").Value = 5.296029437114068
$ws.Range("F15").Value = 35.49596099532999
$ws.Range("G15").Value = 3.614003068182904
$ws.Range("I15").Value = 25.87418645070265
$ws.Range("M15").Value = 22.72391426155774
$ws.Range("B16").Value = 13.67624622392831
$ws.Range("C16").Value = 11.56734666078624
$ws.Range("D16").Value = 5.293023857029421
$ws.Range("F16").Value = 35.01562260650913
$ws.Range("G16").Value = 3.616518934862659
$ws.Range("I16").Value = 25.6706386686088
$ws.Range("M16").Value = 22.27861687699651
$ws.Range("B17").Value = 13.4855328977761
$ws.Range("C17").Value = 11.34396974418851
$ws.Range("D17").Value = 5.291325330141154
$ws.Range("F17").Value = 34.72024924208548
$ws.Range("G17").Value = 3.618093918641153
$ws.Range("I17").Value = 25.54657548083571
$ws.Range("M17").Value = 22.00336013631491
$ws.Range("B18").Value = 13.37497766965713
$ws.Range("C18").Value = 11.21649192887282
$ws.Range("D18").Value = 5.290401303366228
$ws.Range("F18").Value = 34.55016014990085
$ws.Range("G18").Value = 3.61901145511063
$ws.Range("I18").Value = 25.47553166802438
$ws.Range("M18").Value = 21.84433202672322
$ws.Range("B19").Value = 13.33740156314732
$ws.Range("C19").Value = 11.17303707440226
$ws.Range("D19").Value = 5.290097487521658
$ws.Range("F19").Value = 34.49254261820725
$ws.Range("G19").Value = 3.619324121590498
$ws.Range("I19").Value = 25.45153304833912
$ws.Range("M19").Value = 21.7903719078998
$ws.Range("B20").Value = 13.50592481553704
$ws.Range("C20").Value = 11.36742337467497
$ws.Range("D20").Value = 5.291500653322976
$ws.Range("F20").Value = 34.75171410812464
$ws.Range("G20").Value = 3.617925054398422
$ws.Range("I20").Value = 25.55975011021554
$ws.Range("M20").Value = 22.03273625452838
$ws.Range("B21").Value = 14.0615579191735
$ws.Range("C21").Value = 12.04665199011344
$ws.Range("D21").Value = 5.296850724759552
$ws.Range("F21").Value = 35.61971039306346
$ws.Range("G21").Value = 3.613363250611613
$ws.Range("I21").Value = 25.92697352436382
$ws.Range("M21").Value = 22.83819307028009
$ws.Range("B22").Value = 14.41619335764533
$ws.Range("C22").Value = 12.47994308546594
$ws.Range("D22").Value = 5.300828605157418
$ws.Range("F22").Value = 36.18398656541329
$ws.Range("G22").Value = 3.610484685120364
$ws.Range("I22").Value = 26.16939730334546
$ws.Range("M22").Value = 23.35713158049209
$ws.Range("B23").Value = 14.22773284630506
$ws.Range("C23").Value = 12.25055843486618
$ws.Range("D23").Value = 5.298661015938275
$ws.Range("F23").Value = 35.88315565453865
$ws.Range("G23").Value = 3.612011662827915
$ws.Range("I23").Value = 26.03980778334833
$ws.Range("M23").Value = 23.08090154097619
$ws.Range("B24").Value = 13.49670845935658
$ws.Range("C24").Value = 11.35682550961906
$ws.Range("D24").Value = 5.291421226364672
$ws.Range("F24").Value = 34.73748968473009
$ws.Range("G24").Value = 3.618001360356856
$ws.Range("I24").Value = 25.55379297848965
$ws.Range("M24").Value = 22.01945772549675
$ws.Range("B25").Value = 12.67587086397354
$ws.Range("C25").Value = 10.39626743652018
$ws.Range("D25").Value = 5.285628911248629
$ws.Range("F25").Value = 33.49494030333102
$ws.Range("G25").Value = 3.624908063361131
$ws.Range("I25").Value = 25.04176553246887
$ws.Range("M25").Value = 20.84832377353973
